{"js": "// Insert two new achievement bullet paragraphs into the \"KEY ACHIEVEMENTS AND IMPACT\"\n// section, right after the \"Expert methodology validated at highest judicial level\"\n// bullet and before the \"TECHNICAL SKILLS\" heading.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its exact text.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\u2022 Expert methodology validated at highest judicial level\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph for insertion.\");\n}\n\n// First new paragraph: plain bullet text, inserted directly after the anchor.\nconst p1 = anchor.insertParagraph(\n  \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  \"After\"\n);\n\n// Second new paragraph: bullet with a bolded, colored \"178%\" run in the middle.\n// Insert an empty paragraph after p1, then build it up with three runs.\nconst p2 = p1.insertParagraph(\"\", \"After\");\np2.insertText(\"\u2022 \", \"End\");\nconst boldRange = p2.insertText(\"178%\", \"End\");\nboldRange.font.bold = true;\nboldRange.font.color = \"#2C3E50\";\np2.insertText(\" accuracy improvement in racial classification algorithms\", \"End\");\n\nawait context.sync();\n", "ps1": "# Insert two new achievement bullet paragraphs into the \"KEY ACHIEVEMENTS AND IMPACT\"\n# section, right after the \"Expert methodology validated at highest judicial level\"\n# bullet and before the \"TECHNICAL SKILLS\" heading:\n#\n#   \u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\n#   \u2022 178% accuracy improvement in racial classification algorithms   (\"178%\" bold/colored)\n\n$d = $word.ActiveDocument\n$bullet = [char]0x2022\n\n# Locate the anchor paragraph by its exact text (scanning is robust against stale\n# Range/Find handles across structural edits).\n$needle = \"$bullet Expert methodology validated at highest judicial level\"\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\")\n    if ($txt -eq $needle) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph ('Expert methodology validated...') for insertion.\"\n}\n\n# Insert a new empty paragraph right after the anchor. Calling InsertParagraphAfter on the\n# anchor paragraph's OWN Range (collapsed to its end) makes the new paragraph inherit the\n# anchor's (plain/\"Normal\") formatting rather than the following paragraph's.\n$anchorRange = $d.Paragraphs.Item($anchorIndex).Range\n$anchorRange.Collapse(0)   # wdCollapseEnd\n$anchorRange.InsertParagraphAfter()\n\n# --- First new paragraph: plain bullet text ---\n$para1Index = $anchorIndex + 1\n$para1Start = $d.Paragraphs.Item($para1Index).Range.Start\n$ins1 = $d.Range($para1Start, $para1Start)\n$ins1.InsertAfter(\"$bullet Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\")\n\n# Insert a second new empty paragraph right after the first new paragraph, again via that\n# paragraph's own Range so it keeps plain/\"Normal\" formatting (and not the Heading 2 style\n# of the \"TECHNICAL SKILLS\" paragraph that currently follows it).\n$p1Range = $d.Paragraphs.Item($para1Index).Range\n$p1Range.Collapse(0)\n$p1Range.InsertParagraphAfter()\n\n# --- Second new paragraph: three runs - \"\u2022 \", bold/colored \"178%\", then trailing text ---\n$para2Index = $para1Index + 1\n$para2Start = $d.Paragraphs.Item($para2Index).Range.Start\n\n$run1 = $d.Range($para2Start, $para2Start)\n$run1.InsertAfter(\"$bullet \")\n\n$pos2 = $run1.End\n$run2 = $d.Range($pos2, $pos2)\n$run2.InsertAfter(\"178%\")\n# Word/OLE font colors are packed as 0x00BBGGRR (blue, green, red byte order), so\n# RGB(2C,3E,50) becomes 0x503E2C.\n$rVal = 0x2C\n$gVal = 0x3E\n$bVal = 0x50\n$run2.Font.Color = ($bVal * 65536) + ($gVal * 256) + $rVal\n$run2.Font.Bold = 1\n\n$pos3 = $run2.End\n$run3 = $d.Range($pos3, $pos3)\n$run3.InsertAfter(\" accuracy improvement in racial classification algorithms\")\n"}
